# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 14:29"

# Re-rank a block of provincias (city name + stats move together as rows shuffle)
$ws.Range("A11").Value = "Ciudad Real"
$ws.Range("B11").Value = 1543
$ws.Range("C11").Value = 195
$ws.Range("D11").Value = 1356
$ws.Range("E11").Value = 117

$ws.Range("A12").Value = "La Rioja"
$ws.Range("B12").Value = 1436
$ws.Range("C12").Value = 364
$ws.Range("D12").Value = 1007
$ws.Range("E12").Value = 65

$ws.Range("A15").Value = "Toledo"
$ws.Range("B15").Value = 1192
$ws.Range("C15").Value = 195
$ws.Range("D15").Value = 1034
$ws.Range("E15").Value = 98

$ws.Range("A16").Value = "A Coruña"
$ws.Range("B16").Value = 1177
$ws.Range("C16").Value = 95
$ws.Range("D16").Value = 1114
$ws.Range("E16").Value = 31

$ws.Range("A17").Value = "Albacete"
$ws.Range("B17").Value = 1114
$ws.Range("C17").Value = 195
$ws.Range("D17").Value = 968
$ws.Range("E17").Value = 101

$ws.Range("A18").Value = "Malaga"
$ws.Range("B18").Value = 1053
$ws.Range("C18").Value = 80
$ws.Range("D18").Value = 917
$ws.Range("E18").Value = 56

$ws.Range("A19").Value = "Gran Canaria"
$ws.Range("B19").Value = 1025
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 964
$ws.Range("E19").Value = 36

$ws.Range("A20").Value = "Asturias"
$ws.Range("B20").Value = 1004
$ws.Range("C20").Value = 65
$ws.Range("D20").Value = 906
$ws.Range("E20").Value = 33

$ws.Range("A21").Value = "Pontevedra"
$ws.Range("B21").Value = 960
$ws.Range("C21").Value = 95
$ws.Range("D21").Value = 923
$ws.Range("E21").Value = 9

$ws.Range("A22").Value = "Caceres"
$ws.Range("B22").Value = 957
$ws.Range("C22").Value = 38
$ws.Range("D22").Value = 841
$ws.Range("E22").Value = 78

$ws.Range("A23").Value = "Cantabria"
$ws.Range("B23").Value = 937
$ws.Range("C23").Value = 21
$ws.Range("D23").Value = 894
$ws.Range("E23").Value = 22

$ws.Range("A24").Value = "Aragon"
$ws.Range("B24").Value = 907
$ws.Range("C24").Value = 29
$ws.Range("D24").Value = 838
$ws.Range("E24").Value = 40

$ws.Range("A25").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B25").Value = 894
$ws.Range("C25").Value = 1023
$ws.Range("D25").Value = 662
$ws.Range("E25").Value = 29

$ws.Range("A26").Value = "Salamanca"
$ws.Range("B26").Value = 882
$ws.Range("C26").Value = 131
$ws.Range("D26").Value = 667
$ws.Range("E26").Value = 84

$ws.Range("A27").Value = "Sevilla"
$ws.Range("B27").Value = 830
$ws.Range("C27").Value = 13
$ws.Range("D27").Value = 791
$ws.Range("E27").Value = 26

$ws.Range("A28").Value = "Valladolid"
$ws.Range("B28").Value = 807
$ws.Range("C28").Value = 114
$ws.Range("D28").Value = 648
$ws.Range("E28").Value = 45

$ws.Range("A29").Value = "Granada"
$ws.Range("B29").Value = 806
$ws.Range("C29").Value = 11
$ws.Range("D29").Value = 746
$ws.Range("E29").Value = 49

$ws.Range("A30").Value = "Murcia"
$ws.Range("B30").Value = 802
$ws.Range("C30").Value = 12
$ws.Range("D30").Value = 773
$ws.Range("E30").Value = 17

# Guadalajara stats refresh (no reorder)
$ws.Range("B37").Value = 441
$ws.Range("C37").Value = 195
$ws.Range("D37").Value = 355
$ws.Range("E37").Value = 79

# Another small re-rank block
$ws.Range("A45").Value = "Cuenca"
$ws.Range("B45").Value = 222
$ws.Range("C45").Value = 195
$ws.Range("D45").Value = 154
$ws.Range("E45").Value = 53

$ws.Range("A46").Value = "Palencia"
$ws.Range("B46").Value = 220
$ws.Range("C46").Value = 26
$ws.Range("D46").Value = 183
$ws.Range("E46").Value = 11

$ws.Range("A47").Value = "Mallorca"
$ws.Range("B47").Value = 210
$ws.Range("C47").Value = 18
$ws.Range("D47").Value = 194
$ws.Range("E47").Value = 12
